# Ticket 79 - Fix implicit sheet cloning case when the number of items in
# the collection is 1.
#
# The template previously only demonstrated implicit sheet cloning for a
# collection with more than one item (the "${dvs.name}$@l=0" sheet covers
# the l=0 / first division). To exercise (and prove the fix for) the
# single-item-collection case, add a second clone-source sheet,
# "${dvs.name}$@l=1", that is an exact duplicate of "${dvs.name}$@l=0",
# placed right after it in tab order (after "Static3").

$wb = $excel.ActiveWorkbook

# Locate the existing division-template sheet that will be duplicated.
$sourceSheetName = '${dvs.name}$@l=0'
$newSheetName    = '${dvs.name}$@l=1'

$source = $wb.Worksheets.Item($sourceSheetName)

# Remember which sheet is currently active/selected so we can restore the
# original tab selection after we're done (copying/activating a sheet
# moves Excel's active-tab pointer).
$originalActiveSheetName = $wb.ActiveSheet.Name

# Copy the source sheet to the very end of the workbook (after "Static3",
# which is currently the last sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$source.Copy($null, $lastSheet)

# The newly created sheet is now the last sheet in the workbook; rename it
# to match the l=1 division placeholder.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = $newSheetName

# Restore the originally-active sheet/tab.
$wb.Worksheets.Item($originalActiveSheetName).Activate()
